$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cycle: [3, 5]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A3:B3").Value2
$ws.Range("A3:B3").Value2 = $ws.Range("A5:B5").Value2
$ws.Range("A5:B5").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D3:H3").Value2
$ws.Range("D3:H3").Value2 = $ws.Range("D5:H5").Value2
$ws.Range("D5:H5").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M3:M3").Value2
$ws.Range("M3:M3").Value2 = $ws.Range("M5:M5").Value2
$ws.Range("M5:M5").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P3:W3").Value2
$ws.Range("P3:W3").Value2 = $ws.Range("P5:W5").Value2
$ws.Range("P5:W5").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z3:Z3").Value2
$ws.Range("Z3:Z3").Value2 = $ws.Range("Z5:Z5").Value2
$ws.Range("Z5:Z5").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB3:AE3").Value2
$ws.Range("AB3:AE3").Value2 = $ws.Range("AB5:AE5").Value2
$ws.Range("AB5:AE5").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG3:AG3").Value2
$ws.Range("AG3:AG3").Value2 = $ws.Range("AG5:AG5").Value2
$ws.Range("AG5:AG5").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW3:AX3").Value2
$ws.Range("AW3:AX3").Value2 = $ws.Range("AW5:AX5").Value2
$ws.Range("AW5:AX5").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [4, 6]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A4:B4").Value2
$ws.Range("A4:B4").Value2 = $ws.Range("A6:B6").Value2
$ws.Range("A6:B6").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D4:H4").Value2
$ws.Range("D4:H4").Value2 = $ws.Range("D6:H6").Value2
$ws.Range("D6:H6").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M4:M4").Value2
$ws.Range("M4:M4").Value2 = $ws.Range("M6:M6").Value2
$ws.Range("M6:M6").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P4:W4").Value2
$ws.Range("P4:W4").Value2 = $ws.Range("P6:W6").Value2
$ws.Range("P6:W6").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z4:Z4").Value2
$ws.Range("Z4:Z4").Value2 = $ws.Range("Z6:Z6").Value2
$ws.Range("Z6:Z6").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB4:AE4").Value2
$ws.Range("AB4:AE4").Value2 = $ws.Range("AB6:AE6").Value2
$ws.Range("AB6:AE6").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG4:AG4").Value2
$ws.Range("AG4:AG4").Value2 = $ws.Range("AG6:AG6").Value2
$ws.Range("AG6:AG6").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW4:AX4").Value2
$ws.Range("AW4:AX4").Value2 = $ws.Range("AW6:AX6").Value2
$ws.Range("AW6:AX6").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [12, 13]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A12:B12").Value2
$ws.Range("A12:B12").Value2 = $ws.Range("A13:B13").Value2
$ws.Range("A13:B13").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D12:H12").Value2
$ws.Range("D12:H12").Value2 = $ws.Range("D13:H13").Value2
$ws.Range("D13:H13").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M12:M12").Value2
$ws.Range("M12:M12").Value2 = $ws.Range("M13:M13").Value2
$ws.Range("M13:M13").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P12:W12").Value2
$ws.Range("P12:W12").Value2 = $ws.Range("P13:W13").Value2
$ws.Range("P13:W13").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z12:Z12").Value2
$ws.Range("Z12:Z12").Value2 = $ws.Range("Z13:Z13").Value2
$ws.Range("Z13:Z13").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB12:AE12").Value2
$ws.Range("AB12:AE12").Value2 = $ws.Range("AB13:AE13").Value2
$ws.Range("AB13:AE13").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG12:AG12").Value2
$ws.Range("AG12:AG12").Value2 = $ws.Range("AG13:AG13").Value2
$ws.Range("AG13:AG13").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW12:AX12").Value2
$ws.Range("AW12:AX12").Value2 = $ws.Range("AW13:AX13").Value2
$ws.Range("AW13:AX13").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [15, 16, 17, 18, 20]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A15:B15").Value2
$ws.Range("A15:B15").Value2 = $ws.Range("A16:B16").Value2
$ws.Range("A16:B16").Value2 = $ws.Range("A17:B17").Value2
$ws.Range("A17:B17").Value2 = $ws.Range("A18:B18").Value2
$ws.Range("A18:B18").Value2 = $ws.Range("A20:B20").Value2
$ws.Range("A20:B20").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D15:H15").Value2
$ws.Range("D15:H15").Value2 = $ws.Range("D16:H16").Value2
$ws.Range("D16:H16").Value2 = $ws.Range("D17:H17").Value2
$ws.Range("D17:H17").Value2 = $ws.Range("D18:H18").Value2
$ws.Range("D18:H18").Value2 = $ws.Range("D20:H20").Value2
$ws.Range("D20:H20").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M15:M15").Value2
$ws.Range("M15:M15").Value2 = $ws.Range("M16:M16").Value2
$ws.Range("M16:M16").Value2 = $ws.Range("M17:M17").Value2
$ws.Range("M17:M17").Value2 = $ws.Range("M18:M18").Value2
$ws.Range("M18:M18").Value2 = $ws.Range("M20:M20").Value2
$ws.Range("M20:M20").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P15:W15").Value2
$ws.Range("P15:W15").Value2 = $ws.Range("P16:W16").Value2
$ws.Range("P16:W16").Value2 = $ws.Range("P17:W17").Value2
$ws.Range("P17:W17").Value2 = $ws.Range("P18:W18").Value2
$ws.Range("P18:W18").Value2 = $ws.Range("P20:W20").Value2
$ws.Range("P20:W20").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z15:Z15").Value2
$ws.Range("Z15:Z15").Value2 = $ws.Range("Z16:Z16").Value2
$ws.Range("Z16:Z16").Value2 = $ws.Range("Z17:Z17").Value2
$ws.Range("Z17:Z17").Value2 = $ws.Range("Z18:Z18").Value2
$ws.Range("Z18:Z18").Value2 = $ws.Range("Z20:Z20").Value2
$ws.Range("Z20:Z20").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB15:AE15").Value2
$ws.Range("AB15:AE15").Value2 = $ws.Range("AB16:AE16").Value2
$ws.Range("AB16:AE16").Value2 = $ws.Range("AB17:AE17").Value2
$ws.Range("AB17:AE17").Value2 = $ws.Range("AB18:AE18").Value2
$ws.Range("AB18:AE18").Value2 = $ws.Range("AB20:AE20").Value2
$ws.Range("AB20:AE20").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG15:AG15").Value2
$ws.Range("AG15:AG15").Value2 = $ws.Range("AG16:AG16").Value2
$ws.Range("AG16:AG16").Value2 = $ws.Range("AG17:AG17").Value2
$ws.Range("AG17:AG17").Value2 = $ws.Range("AG18:AG18").Value2
$ws.Range("AG18:AG18").Value2 = $ws.Range("AG20:AG20").Value2
$ws.Range("AG20:AG20").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW15:AX15").Value2
$ws.Range("AW15:AX15").Value2 = $ws.Range("AW16:AX16").Value2
$ws.Range("AW16:AX16").Value2 = $ws.Range("AW17:AX17").Value2
$ws.Range("AW17:AX17").Value2 = $ws.Range("AW18:AX18").Value2
$ws.Range("AW18:AX18").Value2 = $ws.Range("AW20:AX20").Value2
$ws.Range("AW20:AX20").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [21, 23, 22, 27, 26]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A21:B21").Value2
$ws.Range("A21:B21").Value2 = $ws.Range("A23:B23").Value2
$ws.Range("A23:B23").Value2 = $ws.Range("A22:B22").Value2
$ws.Range("A22:B22").Value2 = $ws.Range("A27:B27").Value2
$ws.Range("A27:B27").Value2 = $ws.Range("A26:B26").Value2
$ws.Range("A26:B26").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D21:H21").Value2
$ws.Range("D21:H21").Value2 = $ws.Range("D23:H23").Value2
$ws.Range("D23:H23").Value2 = $ws.Range("D22:H22").Value2
$ws.Range("D22:H22").Value2 = $ws.Range("D27:H27").Value2
$ws.Range("D27:H27").Value2 = $ws.Range("D26:H26").Value2
$ws.Range("D26:H26").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M21:M21").Value2
$ws.Range("M21:M21").Value2 = $ws.Range("M23:M23").Value2
$ws.Range("M23:M23").Value2 = $ws.Range("M22:M22").Value2
$ws.Range("M22:M22").Value2 = $ws.Range("M27:M27").Value2
$ws.Range("M27:M27").Value2 = $ws.Range("M26:M26").Value2
$ws.Range("M26:M26").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P21:W21").Value2
$ws.Range("P21:W21").Value2 = $ws.Range("P23:W23").Value2
$ws.Range("P23:W23").Value2 = $ws.Range("P22:W22").Value2
$ws.Range("P22:W22").Value2 = $ws.Range("P27:W27").Value2
$ws.Range("P27:W27").Value2 = $ws.Range("P26:W26").Value2
$ws.Range("P26:W26").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z21:Z21").Value2
$ws.Range("Z21:Z21").Value2 = $ws.Range("Z23:Z23").Value2
$ws.Range("Z23:Z23").Value2 = $ws.Range("Z22:Z22").Value2
$ws.Range("Z22:Z22").Value2 = $ws.Range("Z27:Z27").Value2
$ws.Range("Z27:Z27").Value2 = $ws.Range("Z26:Z26").Value2
$ws.Range("Z26:Z26").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB21:AE21").Value2
$ws.Range("AB21:AE21").Value2 = $ws.Range("AB23:AE23").Value2
$ws.Range("AB23:AE23").Value2 = $ws.Range("AB22:AE22").Value2
$ws.Range("AB22:AE22").Value2 = $ws.Range("AB27:AE27").Value2
$ws.Range("AB27:AE27").Value2 = $ws.Range("AB26:AE26").Value2
$ws.Range("AB26:AE26").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG21:AG21").Value2
$ws.Range("AG21:AG21").Value2 = $ws.Range("AG23:AG23").Value2
$ws.Range("AG23:AG23").Value2 = $ws.Range("AG22:AG22").Value2
$ws.Range("AG22:AG22").Value2 = $ws.Range("AG27:AG27").Value2
$ws.Range("AG27:AG27").Value2 = $ws.Range("AG26:AG26").Value2
$ws.Range("AG26:AG26").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW21:AX21").Value2
$ws.Range("AW21:AX21").Value2 = $ws.Range("AW23:AX23").Value2
$ws.Range("AW23:AX23").Value2 = $ws.Range("AW22:AX22").Value2
$ws.Range("AW22:AX22").Value2 = $ws.Range("AW27:AX27").Value2
$ws.Range("AW27:AX27").Value2 = $ws.Range("AW26:AX26").Value2
$ws.Range("AW26:AX26").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [24, 25]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A24:B24").Value2
$ws.Range("A24:B24").Value2 = $ws.Range("A25:B25").Value2
$ws.Range("A25:B25").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D24:H24").Value2
$ws.Range("D24:H24").Value2 = $ws.Range("D25:H25").Value2
$ws.Range("D25:H25").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M24:M24").Value2
$ws.Range("M24:M24").Value2 = $ws.Range("M25:M25").Value2
$ws.Range("M25:M25").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P24:W24").Value2
$ws.Range("P24:W24").Value2 = $ws.Range("P25:W25").Value2
$ws.Range("P25:W25").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z24:Z24").Value2
$ws.Range("Z24:Z24").Value2 = $ws.Range("Z25:Z25").Value2
$ws.Range("Z25:Z25").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB24:AE24").Value2
$ws.Range("AB24:AE24").Value2 = $ws.Range("AB25:AE25").Value2
$ws.Range("AB25:AE25").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG24:AG24").Value2
$ws.Range("AG24:AG24").Value2 = $ws.Range("AG25:AG25").Value2
$ws.Range("AG25:AG25").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW24:AX24").Value2
$ws.Range("AW24:AX24").Value2 = $ws.Range("AW25:AX25").Value2
$ws.Range("AW25:AX25").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [34, 35]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A34:B34").Value2
$ws.Range("A34:B34").Value2 = $ws.Range("A35:B35").Value2
$ws.Range("A35:B35").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D34:H34").Value2
$ws.Range("D34:H34").Value2 = $ws.Range("D35:H35").Value2
$ws.Range("D35:H35").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M34:M34").Value2
$ws.Range("M34:M34").Value2 = $ws.Range("M35:M35").Value2
$ws.Range("M35:M35").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P34:W34").Value2
$ws.Range("P34:W34").Value2 = $ws.Range("P35:W35").Value2
$ws.Range("P35:W35").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z34:Z34").Value2
$ws.Range("Z34:Z34").Value2 = $ws.Range("Z35:Z35").Value2
$ws.Range("Z35:Z35").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB34:AE34").Value2
$ws.Range("AB34:AE34").Value2 = $ws.Range("AB35:AE35").Value2
$ws.Range("AB35:AE35").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG34:AG34").Value2
$ws.Range("AG34:AG34").Value2 = $ws.Range("AG35:AG35").Value2
$ws.Range("AG35:AG35").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW34:AX34").Value2
$ws.Range("AW34:AX34").Value2 = $ws.Range("AW35:AX35").Value2
$ws.Range("AW35:AX35").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [50, 51, 52, 53]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A50:B50").Value2
$ws.Range("A50:B50").Value2 = $ws.Range("A51:B51").Value2
$ws.Range("A51:B51").Value2 = $ws.Range("A52:B52").Value2
$ws.Range("A52:B52").Value2 = $ws.Range("A53:B53").Value2
$ws.Range("A53:B53").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D50:H50").Value2
$ws.Range("D50:H50").Value2 = $ws.Range("D51:H51").Value2
$ws.Range("D51:H51").Value2 = $ws.Range("D52:H52").Value2
$ws.Range("D52:H52").Value2 = $ws.Range("D53:H53").Value2
$ws.Range("D53:H53").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M50:M50").Value2
$ws.Range("M50:M50").Value2 = $ws.Range("M51:M51").Value2
$ws.Range("M51:M51").Value2 = $ws.Range("M52:M52").Value2
$ws.Range("M52:M52").Value2 = $ws.Range("M53:M53").Value2
$ws.Range("M53:M53").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P50:W50").Value2
$ws.Range("P50:W50").Value2 = $ws.Range("P51:W51").Value2
$ws.Range("P51:W51").Value2 = $ws.Range("P52:W52").Value2
$ws.Range("P52:W52").Value2 = $ws.Range("P53:W53").Value2
$ws.Range("P53:W53").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z50:Z50").Value2
$ws.Range("Z50:Z50").Value2 = $ws.Range("Z51:Z51").Value2
$ws.Range("Z51:Z51").Value2 = $ws.Range("Z52:Z52").Value2
$ws.Range("Z52:Z52").Value2 = $ws.Range("Z53:Z53").Value2
$ws.Range("Z53:Z53").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB50:AE50").Value2
$ws.Range("AB50:AE50").Value2 = $ws.Range("AB51:AE51").Value2
$ws.Range("AB51:AE51").Value2 = $ws.Range("AB52:AE52").Value2
$ws.Range("AB52:AE52").Value2 = $ws.Range("AB53:AE53").Value2
$ws.Range("AB53:AE53").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG50:AG50").Value2
$ws.Range("AG50:AG50").Value2 = $ws.Range("AG51:AG51").Value2
$ws.Range("AG51:AG51").Value2 = $ws.Range("AG52:AG52").Value2
$ws.Range("AG52:AG52").Value2 = $ws.Range("AG53:AG53").Value2
$ws.Range("AG53:AG53").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW50:AX50").Value2
$ws.Range("AW50:AX50").Value2 = $ws.Range("AW51:AX51").Value2
$ws.Range("AW51:AX51").Value2 = $ws.Range("AW52:AX52").Value2
$ws.Range("AW52:AX52").Value2 = $ws.Range("AW53:AX53").Value2
$ws.Range("AW53:AX53").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [55, 60, 59, 58, 57, 56]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A55:B55").Value2
$ws.Range("A55:B55").Value2 = $ws.Range("A60:B60").Value2
$ws.Range("A60:B60").Value2 = $ws.Range("A59:B59").Value2
$ws.Range("A59:B59").Value2 = $ws.Range("A58:B58").Value2
$ws.Range("A58:B58").Value2 = $ws.Range("A57:B57").Value2
$ws.Range("A57:B57").Value2 = $ws.Range("A56:B56").Value2
$ws.Range("A56:B56").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D55:H55").Value2
$ws.Range("D55:H55").Value2 = $ws.Range("D60:H60").Value2
$ws.Range("D60:H60").Value2 = $ws.Range("D59:H59").Value2
$ws.Range("D59:H59").Value2 = $ws.Range("D58:H58").Value2
$ws.Range("D58:H58").Value2 = $ws.Range("D57:H57").Value2
$ws.Range("D57:H57").Value2 = $ws.Range("D56:H56").Value2
$ws.Range("D56:H56").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M55:M55").Value2
$ws.Range("M55:M55").Value2 = $ws.Range("M60:M60").Value2
$ws.Range("M60:M60").Value2 = $ws.Range("M59:M59").Value2
$ws.Range("M59:M59").Value2 = $ws.Range("M58:M58").Value2
$ws.Range("M58:M58").Value2 = $ws.Range("M57:M57").Value2
$ws.Range("M57:M57").Value2 = $ws.Range("M56:M56").Value2
$ws.Range("M56:M56").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P55:W55").Value2
$ws.Range("P55:W55").Value2 = $ws.Range("P60:W60").Value2
$ws.Range("P60:W60").Value2 = $ws.Range("P59:W59").Value2
$ws.Range("P59:W59").Value2 = $ws.Range("P58:W58").Value2
$ws.Range("P58:W58").Value2 = $ws.Range("P57:W57").Value2
$ws.Range("P57:W57").Value2 = $ws.Range("P56:W56").Value2
$ws.Range("P56:W56").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z55:Z55").Value2
$ws.Range("Z55:Z55").Value2 = $ws.Range("Z60:Z60").Value2
$ws.Range("Z60:Z60").Value2 = $ws.Range("Z59:Z59").Value2
$ws.Range("Z59:Z59").Value2 = $ws.Range("Z58:Z58").Value2
$ws.Range("Z58:Z58").Value2 = $ws.Range("Z57:Z57").Value2
$ws.Range("Z57:Z57").Value2 = $ws.Range("Z56:Z56").Value2
$ws.Range("Z56:Z56").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB55:AE55").Value2
$ws.Range("AB55:AE55").Value2 = $ws.Range("AB60:AE60").Value2
$ws.Range("AB60:AE60").Value2 = $ws.Range("AB59:AE59").Value2
$ws.Range("AB59:AE59").Value2 = $ws.Range("AB58:AE58").Value2
$ws.Range("AB58:AE58").Value2 = $ws.Range("AB57:AE57").Value2
$ws.Range("AB57:AE57").Value2 = $ws.Range("AB56:AE56").Value2
$ws.Range("AB56:AE56").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG55:AG55").Value2
$ws.Range("AG55:AG55").Value2 = $ws.Range("AG60:AG60").Value2
$ws.Range("AG60:AG60").Value2 = $ws.Range("AG59:AG59").Value2
$ws.Range("AG59:AG59").Value2 = $ws.Range("AG58:AG58").Value2
$ws.Range("AG58:AG58").Value2 = $ws.Range("AG57:AG57").Value2
$ws.Range("AG57:AG57").Value2 = $ws.Range("AG56:AG56").Value2
$ws.Range("AG56:AG56").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW55:AX55").Value2
$ws.Range("AW55:AX55").Value2 = $ws.Range("AW60:AX60").Value2
$ws.Range("AW60:AX60").Value2 = $ws.Range("AW59:AX59").Value2
$ws.Range("AW59:AX59").Value2 = $ws.Range("AW58:AX58").Value2
$ws.Range("AW58:AX58").Value2 = $ws.Range("AW57:AX57").Value2
$ws.Range("AW57:AX57").Value2 = $ws.Range("AW56:AX56").Value2
$ws.Range("AW56:AX56").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [62, 64, 63]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A62:B62").Value2
$ws.Range("A62:B62").Value2 = $ws.Range("A64:B64").Value2
$ws.Range("A64:B64").Value2 = $ws.Range("A63:B63").Value2
$ws.Range("A63:B63").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D62:H62").Value2
$ws.Range("D62:H62").Value2 = $ws.Range("D64:H64").Value2
$ws.Range("D64:H64").Value2 = $ws.Range("D63:H63").Value2
$ws.Range("D63:H63").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M62:M62").Value2
$ws.Range("M62:M62").Value2 = $ws.Range("M64:M64").Value2
$ws.Range("M64:M64").Value2 = $ws.Range("M63:M63").Value2
$ws.Range("M63:M63").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P62:W62").Value2
$ws.Range("P62:W62").Value2 = $ws.Range("P64:W64").Value2
$ws.Range("P64:W64").Value2 = $ws.Range("P63:W63").Value2
$ws.Range("P63:W63").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z62:Z62").Value2
$ws.Range("Z62:Z62").Value2 = $ws.Range("Z64:Z64").Value2
$ws.Range("Z64:Z64").Value2 = $ws.Range("Z63:Z63").Value2
$ws.Range("Z63:Z63").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB62:AE62").Value2
$ws.Range("AB62:AE62").Value2 = $ws.Range("AB64:AE64").Value2
$ws.Range("AB64:AE64").Value2 = $ws.Range("AB63:AE63").Value2
$ws.Range("AB63:AE63").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG62:AG62").Value2
$ws.Range("AG62:AG62").Value2 = $ws.Range("AG64:AG64").Value2
$ws.Range("AG64:AG64").Value2 = $ws.Range("AG63:AG63").Value2
$ws.Range("AG63:AG63").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW62:AX62").Value2
$ws.Range("AW62:AX62").Value2 = $ws.Range("AW64:AX64").Value2
$ws.Range("AW64:AX64").Value2 = $ws.Range("AW63:AX63").Value2
$ws.Range("AW63:AX63").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [65, 66]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A65:B65").Value2
$ws.Range("A65:B65").Value2 = $ws.Range("A66:B66").Value2
$ws.Range("A66:B66").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D65:H65").Value2
$ws.Range("D65:H65").Value2 = $ws.Range("D66:H66").Value2
$ws.Range("D66:H66").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M65:M65").Value2
$ws.Range("M65:M65").Value2 = $ws.Range("M66:M66").Value2
$ws.Range("M66:M66").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P65:W65").Value2
$ws.Range("P65:W65").Value2 = $ws.Range("P66:W66").Value2
$ws.Range("P66:W66").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z65:Z65").Value2
$ws.Range("Z65:Z65").Value2 = $ws.Range("Z66:Z66").Value2
$ws.Range("Z66:Z66").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB65:AE65").Value2
$ws.Range("AB65:AE65").Value2 = $ws.Range("AB66:AE66").Value2
$ws.Range("AB66:AE66").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG65:AG65").Value2
$ws.Range("AG65:AG65").Value2 = $ws.Range("AG66:AG66").Value2
$ws.Range("AG66:AG66").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW65:AX65").Value2
$ws.Range("AW65:AX65").Value2 = $ws.Range("AW66:AX66").Value2
$ws.Range("AW66:AX66").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [68, 71, 69, 70]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A68:B68").Value2
$ws.Range("A68:B68").Value2 = $ws.Range("A71:B71").Value2
$ws.Range("A71:B71").Value2 = $ws.Range("A69:B69").Value2
$ws.Range("A69:B69").Value2 = $ws.Range("A70:B70").Value2
$ws.Range("A70:B70").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D68:H68").Value2
$ws.Range("D68:H68").Value2 = $ws.Range("D71:H71").Value2
$ws.Range("D71:H71").Value2 = $ws.Range("D69:H69").Value2
$ws.Range("D69:H69").Value2 = $ws.Range("D70:H70").Value2
$ws.Range("D70:H70").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M68:M68").Value2
$ws.Range("M68:M68").Value2 = $ws.Range("M71:M71").Value2
$ws.Range("M71:M71").Value2 = $ws.Range("M69:M69").Value2
$ws.Range("M69:M69").Value2 = $ws.Range("M70:M70").Value2
$ws.Range("M70:M70").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P68:W68").Value2
$ws.Range("P68:W68").Value2 = $ws.Range("P71:W71").Value2
$ws.Range("P71:W71").Value2 = $ws.Range("P69:W69").Value2
$ws.Range("P69:W69").Value2 = $ws.Range("P70:W70").Value2
$ws.Range("P70:W70").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z68:Z68").Value2
$ws.Range("Z68:Z68").Value2 = $ws.Range("Z71:Z71").Value2
$ws.Range("Z71:Z71").Value2 = $ws.Range("Z69:Z69").Value2
$ws.Range("Z69:Z69").Value2 = $ws.Range("Z70:Z70").Value2
$ws.Range("Z70:Z70").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB68:AE68").Value2
$ws.Range("AB68:AE68").Value2 = $ws.Range("AB71:AE71").Value2
$ws.Range("AB71:AE71").Value2 = $ws.Range("AB69:AE69").Value2
$ws.Range("AB69:AE69").Value2 = $ws.Range("AB70:AE70").Value2
$ws.Range("AB70:AE70").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG68:AG68").Value2
$ws.Range("AG68:AG68").Value2 = $ws.Range("AG71:AG71").Value2
$ws.Range("AG71:AG71").Value2 = $ws.Range("AG69:AG69").Value2
$ws.Range("AG69:AG69").Value2 = $ws.Range("AG70:AG70").Value2
$ws.Range("AG70:AG70").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW68:AX68").Value2
$ws.Range("AW68:AX68").Value2 = $ws.Range("AW71:AX71").Value2
$ws.Range("AW71:AX71").Value2 = $ws.Range("AW69:AX69").Value2
$ws.Range("AW69:AX69").Value2 = $ws.Range("AW70:AX70").Value2
$ws.Range("AW70:AX70").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [72, 73, 74, 75]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A72:B72").Value2
$ws.Range("A72:B72").Value2 = $ws.Range("A73:B73").Value2
$ws.Range("A73:B73").Value2 = $ws.Range("A74:B74").Value2
$ws.Range("A74:B74").Value2 = $ws.Range("A75:B75").Value2
$ws.Range("A75:B75").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D72:H72").Value2
$ws.Range("D72:H72").Value2 = $ws.Range("D73:H73").Value2
$ws.Range("D73:H73").Value2 = $ws.Range("D74:H74").Value2
$ws.Range("D74:H74").Value2 = $ws.Range("D75:H75").Value2
$ws.Range("D75:H75").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M72:M72").Value2
$ws.Range("M72:M72").Value2 = $ws.Range("M73:M73").Value2
$ws.Range("M73:M73").Value2 = $ws.Range("M74:M74").Value2
$ws.Range("M74:M74").Value2 = $ws.Range("M75:M75").Value2
$ws.Range("M75:M75").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P72:W72").Value2
$ws.Range("P72:W72").Value2 = $ws.Range("P73:W73").Value2
$ws.Range("P73:W73").Value2 = $ws.Range("P74:W74").Value2
$ws.Range("P74:W74").Value2 = $ws.Range("P75:W75").Value2
$ws.Range("P75:W75").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z72:Z72").Value2
$ws.Range("Z72:Z72").Value2 = $ws.Range("Z73:Z73").Value2
$ws.Range("Z73:Z73").Value2 = $ws.Range("Z74:Z74").Value2
$ws.Range("Z74:Z74").Value2 = $ws.Range("Z75:Z75").Value2
$ws.Range("Z75:Z75").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB72:AE72").Value2
$ws.Range("AB72:AE72").Value2 = $ws.Range("AB73:AE73").Value2
$ws.Range("AB73:AE73").Value2 = $ws.Range("AB74:AE74").Value2
$ws.Range("AB74:AE74").Value2 = $ws.Range("AB75:AE75").Value2
$ws.Range("AB75:AE75").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG72:AG72").Value2
$ws.Range("AG72:AG72").Value2 = $ws.Range("AG73:AG73").Value2
$ws.Range("AG73:AG73").Value2 = $ws.Range("AG74:AG74").Value2
$ws.Range("AG74:AG74").Value2 = $ws.Range("AG75:AG75").Value2
$ws.Range("AG75:AG75").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW72:AX72").Value2
$ws.Range("AW72:AX72").Value2 = $ws.Range("AW73:AX73").Value2
$ws.Range("AW73:AX73").Value2 = $ws.Range("AW74:AX74").Value2
$ws.Range("AW74:AX74").Value2 = $ws.Range("AW75:AX75").Value2
$ws.Range("AW75:AX75").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [80, 86, 81]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A80:B80").Value2
$ws.Range("A80:B80").Value2 = $ws.Range("A86:B86").Value2
$ws.Range("A86:B86").Value2 = $ws.Range("A81:B81").Value2
$ws.Range("A81:B81").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D80:H80").Value2
$ws.Range("D80:H80").Value2 = $ws.Range("D86:H86").Value2
$ws.Range("D86:H86").Value2 = $ws.Range("D81:H81").Value2
$ws.Range("D81:H81").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M80:M80").Value2
$ws.Range("M80:M80").Value2 = $ws.Range("M86:M86").Value2
$ws.Range("M86:M86").Value2 = $ws.Range("M81:M81").Value2
$ws.Range("M81:M81").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P80:W80").Value2
$ws.Range("P80:W80").Value2 = $ws.Range("P86:W86").Value2
$ws.Range("P86:W86").Value2 = $ws.Range("P81:W81").Value2
$ws.Range("P81:W81").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z80:Z80").Value2
$ws.Range("Z80:Z80").Value2 = $ws.Range("Z86:Z86").Value2
$ws.Range("Z86:Z86").Value2 = $ws.Range("Z81:Z81").Value2
$ws.Range("Z81:Z81").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB80:AE80").Value2
$ws.Range("AB80:AE80").Value2 = $ws.Range("AB86:AE86").Value2
$ws.Range("AB86:AE86").Value2 = $ws.Range("AB81:AE81").Value2
$ws.Range("AB81:AE81").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG80:AG80").Value2
$ws.Range("AG80:AG80").Value2 = $ws.Range("AG86:AG86").Value2
$ws.Range("AG86:AG86").Value2 = $ws.Range("AG81:AG81").Value2
$ws.Range("AG81:AG81").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW80:AX80").Value2
$ws.Range("AW80:AX80").Value2 = $ws.Range("AW86:AX86").Value2
$ws.Range("AW86:AX86").Value2 = $ws.Range("AW81:AX81").Value2
$ws.Range("AW81:AX81").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [88, 89]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A88:B88").Value2
$ws.Range("A88:B88").Value2 = $ws.Range("A89:B89").Value2
$ws.Range("A89:B89").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D88:H88").Value2
$ws.Range("D88:H88").Value2 = $ws.Range("D89:H89").Value2
$ws.Range("D89:H89").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M88:M88").Value2
$ws.Range("M88:M88").Value2 = $ws.Range("M89:M89").Value2
$ws.Range("M89:M89").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P88:W88").Value2
$ws.Range("P88:W88").Value2 = $ws.Range("P89:W89").Value2
$ws.Range("P89:W89").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z88:Z88").Value2
$ws.Range("Z88:Z88").Value2 = $ws.Range("Z89:Z89").Value2
$ws.Range("Z89:Z89").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB88:AE88").Value2
$ws.Range("AB88:AE88").Value2 = $ws.Range("AB89:AE89").Value2
$ws.Range("AB89:AE89").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG88:AG88").Value2
$ws.Range("AG88:AG88").Value2 = $ws.Range("AG89:AG89").Value2
$ws.Range("AG89:AG89").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW88:AX88").Value2
$ws.Range("AW88:AX88").Value2 = $ws.Range("AW89:AX89").Value2
$ws.Range("AW89:AX89").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [92, 95, 93, 94]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A92:B92").Value2
$ws.Range("A92:B92").Value2 = $ws.Range("A95:B95").Value2
$ws.Range("A95:B95").Value2 = $ws.Range("A93:B93").Value2
$ws.Range("A93:B93").Value2 = $ws.Range("A94:B94").Value2
$ws.Range("A94:B94").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D92:H92").Value2
$ws.Range("D92:H92").Value2 = $ws.Range("D95:H95").Value2
$ws.Range("D95:H95").Value2 = $ws.Range("D93:H93").Value2
$ws.Range("D93:H93").Value2 = $ws.Range("D94:H94").Value2
$ws.Range("D94:H94").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M92:M92").Value2
$ws.Range("M92:M92").Value2 = $ws.Range("M95:M95").Value2
$ws.Range("M95:M95").Value2 = $ws.Range("M93:M93").Value2
$ws.Range("M93:M93").Value2 = $ws.Range("M94:M94").Value2
$ws.Range("M94:M94").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P92:W92").Value2
$ws.Range("P92:W92").Value2 = $ws.Range("P95:W95").Value2
$ws.Range("P95:W95").Value2 = $ws.Range("P93:W93").Value2
$ws.Range("P93:W93").Value2 = $ws.Range("P94:W94").Value2
$ws.Range("P94:W94").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z92:Z92").Value2
$ws.Range("Z92:Z92").Value2 = $ws.Range("Z95:Z95").Value2
$ws.Range("Z95:Z95").Value2 = $ws.Range("Z93:Z93").Value2
$ws.Range("Z93:Z93").Value2 = $ws.Range("Z94:Z94").Value2
$ws.Range("Z94:Z94").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB92:AE92").Value2
$ws.Range("AB92:AE92").Value2 = $ws.Range("AB95:AE95").Value2
$ws.Range("AB95:AE95").Value2 = $ws.Range("AB93:AE93").Value2
$ws.Range("AB93:AE93").Value2 = $ws.Range("AB94:AE94").Value2
$ws.Range("AB94:AE94").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG92:AG92").Value2
$ws.Range("AG92:AG92").Value2 = $ws.Range("AG95:AG95").Value2
$ws.Range("AG95:AG95").Value2 = $ws.Range("AG93:AG93").Value2
$ws.Range("AG93:AG93").Value2 = $ws.Range("AG94:AG94").Value2
$ws.Range("AG94:AG94").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW92:AX92").Value2
$ws.Range("AW92:AX92").Value2 = $ws.Range("AW95:AX95").Value2
$ws.Range("AW95:AX95").Value2 = $ws.Range("AW93:AX93").Value2
$ws.Range("AW93:AX93").Value2 = $ws.Range("AW94:AX94").Value2
$ws.Range("AW94:AX94").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [96, 101, 99]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A96:B96").Value2
$ws.Range("A96:B96").Value2 = $ws.Range("A101:B101").Value2
$ws.Range("A101:B101").Value2 = $ws.Range("A99:B99").Value2
$ws.Range("A99:B99").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D96:H96").Value2
$ws.Range("D96:H96").Value2 = $ws.Range("D101:H101").Value2
$ws.Range("D101:H101").Value2 = $ws.Range("D99:H99").Value2
$ws.Range("D99:H99").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M96:M96").Value2
$ws.Range("M96:M96").Value2 = $ws.Range("M101:M101").Value2
$ws.Range("M101:M101").Value2 = $ws.Range("M99:M99").Value2
$ws.Range("M99:M99").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P96:W96").Value2
$ws.Range("P96:W96").Value2 = $ws.Range("P101:W101").Value2
$ws.Range("P101:W101").Value2 = $ws.Range("P99:W99").Value2
$ws.Range("P99:W99").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z96:Z96").Value2
$ws.Range("Z96:Z96").Value2 = $ws.Range("Z101:Z101").Value2
$ws.Range("Z101:Z101").Value2 = $ws.Range("Z99:Z99").Value2
$ws.Range("Z99:Z99").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB96:AE96").Value2
$ws.Range("AB96:AE96").Value2 = $ws.Range("AB101:AE101").Value2
$ws.Range("AB101:AE101").Value2 = $ws.Range("AB99:AE99").Value2
$ws.Range("AB99:AE99").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG96:AG96").Value2
$ws.Range("AG96:AG96").Value2 = $ws.Range("AG101:AG101").Value2
$ws.Range("AG101:AG101").Value2 = $ws.Range("AG99:AG99").Value2
$ws.Range("AG99:AG99").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW96:AX96").Value2
$ws.Range("AW96:AX96").Value2 = $ws.Range("AW101:AX101").Value2
$ws.Range("AW101:AX101").Value2 = $ws.Range("AW99:AX99").Value2
$ws.Range("AW99:AX99").Value2 = $ws.Range("AW1000:AX1000").Value2

# Cycle: [97, 98, 102, 100]
$ws.Range("A1000:B1000").Value2 = $ws.Range("A97:B97").Value2
$ws.Range("A97:B97").Value2 = $ws.Range("A98:B98").Value2
$ws.Range("A98:B98").Value2 = $ws.Range("A102:B102").Value2
$ws.Range("A102:B102").Value2 = $ws.Range("A100:B100").Value2
$ws.Range("A100:B100").Value2 = $ws.Range("A1000:B1000").Value2
$ws.Range("D1000:H1000").Value2 = $ws.Range("D97:H97").Value2
$ws.Range("D97:H97").Value2 = $ws.Range("D98:H98").Value2
$ws.Range("D98:H98").Value2 = $ws.Range("D102:H102").Value2
$ws.Range("D102:H102").Value2 = $ws.Range("D100:H100").Value2
$ws.Range("D100:H100").Value2 = $ws.Range("D1000:H1000").Value2
$ws.Range("M1000:M1000").Value2 = $ws.Range("M97:M97").Value2
$ws.Range("M97:M97").Value2 = $ws.Range("M98:M98").Value2
$ws.Range("M98:M98").Value2 = $ws.Range("M102:M102").Value2
$ws.Range("M102:M102").Value2 = $ws.Range("M100:M100").Value2
$ws.Range("M100:M100").Value2 = $ws.Range("M1000:M1000").Value2
$ws.Range("P1000:W1000").Value2 = $ws.Range("P97:W97").Value2
$ws.Range("P97:W97").Value2 = $ws.Range("P98:W98").Value2
$ws.Range("P98:W98").Value2 = $ws.Range("P102:W102").Value2
$ws.Range("P102:W102").Value2 = $ws.Range("P100:W100").Value2
$ws.Range("P100:W100").Value2 = $ws.Range("P1000:W1000").Value2
$ws.Range("Z1000:Z1000").Value2 = $ws.Range("Z97:Z97").Value2
$ws.Range("Z97:Z97").Value2 = $ws.Range("Z98:Z98").Value2
$ws.Range("Z98:Z98").Value2 = $ws.Range("Z102:Z102").Value2
$ws.Range("Z102:Z102").Value2 = $ws.Range("Z100:Z100").Value2
$ws.Range("Z100:Z100").Value2 = $ws.Range("Z1000:Z1000").Value2
$ws.Range("AB1000:AE1000").Value2 = $ws.Range("AB97:AE97").Value2
$ws.Range("AB97:AE97").Value2 = $ws.Range("AB98:AE98").Value2
$ws.Range("AB98:AE98").Value2 = $ws.Range("AB102:AE102").Value2
$ws.Range("AB102:AE102").Value2 = $ws.Range("AB100:AE100").Value2
$ws.Range("AB100:AE100").Value2 = $ws.Range("AB1000:AE1000").Value2
$ws.Range("AG1000:AG1000").Value2 = $ws.Range("AG97:AG97").Value2
$ws.Range("AG97:AG97").Value2 = $ws.Range("AG98:AG98").Value2
$ws.Range("AG98:AG98").Value2 = $ws.Range("AG102:AG102").Value2
$ws.Range("AG102:AG102").Value2 = $ws.Range("AG100:AG100").Value2
$ws.Range("AG100:AG100").Value2 = $ws.Range("AG1000:AG1000").Value2
$ws.Range("AW1000:AX1000").Value2 = $ws.Range("AW97:AX97").Value2
$ws.Range("AW97:AX97").Value2 = $ws.Range("AW98:AX98").Value2
$ws.Range("AW98:AX98").Value2 = $ws.Range("AW102:AX102").Value2
$ws.Range("AW102:AX102").Value2 = $ws.Range("AW100:AX100").Value2
$ws.Range("AW100:AX100").Value2 = $ws.Range("AW1000:AX1000").Value2

$ws.Range("A1000:AY1000").ClearContents()